$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "parameters" sheet: insert two new rows (cavalry_distance,
#    cavalry_height_difference_threshold) before the existing
#    "flier_distance" row.
# ---------------------------------------------------------------
$params = $wb.Worksheets.Item("parameters")

$params.Rows("8:9").Insert()

$params.Cells.Item(8, 1).Value = "cavalry_distance"
$params.Cells.Item(8, 2).Value = 4.5

$params.Cells.Item(9, 1).Value = "cavalry_height_difference_threshold"
$params.Cells.Item(9, 2).Value = 2

# Leave this sheet no longer the active tab, matching the recorded
# selection state.
$params.Range("A8:A9").EntireRow.Select()

# ---------------------------------------------------------------
# 2) "interactions" sheet: insert a new "cavalry" interaction
#    column before the existing "siege" column (F), duplicating the
#    per-row values already present in the other columns.
# ---------------------------------------------------------------
$inter = $wb.Worksheets.Item("interactions")

$inter.Columns("F:F").Insert()

$inter.Cells.Item(1, 6).Value = "cavalry"

$inter.Cells.Item(2, 6).Value = -1
$inter.Cells.Item(3, 6).Value = -1
$inter.Cells.Item(4, 6).Value = 0
$inter.Cells.Item(5, 6).Value = 0
$inter.Cells.Item(6, 6).Value = 0
$inter.Cells.Item(7, 6).Value = 0
$inter.Cells.Item(8, 6).Value = 0
$inter.Cells.Item(9, 6).Value = 0
$inter.Cells.Item(10, 6).Value = 0
$inter.Cells.Item(11, 6).Value = 0
$inter.Cells.Item(12, 6).Value = 0
$inter.Cells.Item(13, 6).Value = 0
$inter.Cells.Item(14, 6).Value = -1
$inter.Cells.Item(15, 6).Value = -1
$inter.Cells.Item(16, 6).Value = -1
$inter.Cells.Item(17, 6).Value = -1

# The whole data block (now C:G) loses the bold/italic "header-ish"
# look it inherited from the original siege column formatting.
$dataBlock = $inter.Range("C2:G17")
$dataBlock.Font.Bold = $false
$dataBlock.Font.Italic = $false

# Make "interactions" the active sheet/tab, selection anchored on the
# new header cell.
$inter.Activate()
$inter.Range("F1").Select()
